# Round the ConvexHullArea (column D) values down to whole numbers.
# The workbook's raw/high-precision area measurements are replaced with
# their rounded integer equivalents (values were exported with extra
# floating point precision; this normalizes them to integers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$roundedValues = @{
    2  = 306302
    3  = 152249
    4  = 146107
    5  = 873037
    6  = 1073710
    7  = 11511
    8  = 212456
    9  = 779886
    10 = 419338
    11 = 392631
    12 = 579439
    13 = 375792
    14 = 546909
    15 = 391367
    16 = 834608
    17 = 508669
    18 = 1255156
    19 = 606737
    20 = 62758
    21 = 221311
    22 = 368426
    23 = 525619
    24 = 2571450
    25 = 51948
    26 = 2363908
    27 = 196516
}

foreach ($row in $roundedValues.Keys) {
    $ws.Range("D$row").Value = $roundedValues[$row]
}
